$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
# This string appears on the "Overview" sheet (columns E:F, rows 2-4)
# and on the per-locale sheets "zh-cn" / "de-de" (column C, rows 2-4).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- Narrow the per-locale status columns ---
# Overview sheet: columns E (zh-cn) and F (de-de) shrink from ~17.22 to ~13.41 chars
$wsOverview.Columns.Item(5).ColumnWidth = 12.42
$wsOverview.Columns.Item(6).ColumnWidth = 12.42

# zh-cn / de-de sheets: column C (Status) shrinks from ~17.22 to ~13.41 chars
$wsZhCn.Columns.Item(3).ColumnWidth = 12.42
$wsDeDe.Columns.Item(3).ColumnWidth = 12.42
